# Adds a new slide 5 "Patrón de Diseño (Strategy) / SRP" after the
# existing 4 slides, with a title+content placeholder pair plus a
# second title/content-style textbox pair below it (SRP explanation).

$p = $ppt.ActivePresentation

# New slide, using the same "Título y objetos" (Title and Content)
# layout that the other content slides (2-4) use -> CustomLayout index 2.
$s = $p.Slides.Add(5, 2)

# --- Shape 1: Title placeholder ------------------------------------
$title = $s.Shapes.Item(1)
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Patrón de Diseño (Strategy):"
$titleTr.LanguageID = "es-CO"

# --- Shape 2: Content placeholder (idx=1) ---------------------------
$content = $s.Shapes.Item(2)
$content.Left = 66
$content.Top = 143.75
$content.Width = 828
$content.Height = 96.25
$content.TextFrame.AutoSize = 2
$contentTr = $content.TextFrame.TextRange
$contentTr.Text = "El patrón de diseño Strategy es un patrón de comportamiento que permite definir una familia de algoritmos, encapsular cada uno de ellos y hacerlos intercambiables."
$contentTr.LanguageID = "es-MX"

# --- Shape 3: second "title" textbox (SRP heading) ------------------
$srpTitle = $s.Shapes.AddTextbox(1, 66, 240, 828, 104.37503937007874)
$srpTitle.Name = "Título 1"
$srpTitle.TextFrame.MarginLeft = 7.2
$srpTitle.TextFrame.MarginTop = 3.6
$srpTitle.TextFrame.MarginRight = 7.2
$srpTitle.TextFrame.MarginBottom = 3.6
$srpTitle.TextFrame2.Orientation = 1
$srpTitle.TextFrame2.VerticalAnchor = 3
$srpTitle.TextFrame.AutoSize = 2
$srpTitleTr = $srpTitle.TextFrame.TextRange
$srpTitleTr.Text = "Principio de Responsabilidad Única (SRP):"
$srpTitleTr.LanguageID = "es-CO"
$srpTitleTr.Font.Size = 44

# --- Shape 4: second "content" textbox (SRP explanation) ------------
$srpBody = $s.Shapes.AddTextbox(1, 66, 355, 828, 96.25)
$srpBody.Name = "Marcador de contenido 2"
$srpBody.TextFrame.MarginLeft = 7.2
$srpBody.TextFrame.MarginTop = 3.6
$srpBody.TextFrame.MarginRight = 7.2
$srpBody.TextFrame.MarginBottom = 3.6
$srpBody.TextFrame2.Orientation = 1
$srpBody.TextFrame.AutoSize = 2
$srpBodyTr = $srpBody.TextFrame.TextRange
$srpBodyTr.Text = "El principio de responsabilidad única establece que una clase debe tener solo una razón para cambia.`r"
$srpBodyTr.LanguageID = "es-MX"

$firstPara = $srpBodyTr.Paragraphs(1, 1)
$firstPara.Font.Size = 28

$trailingPara = $srpBodyTr.Paragraphs(2, 1)
$trailingPara.Text = ""
$trailingPara.ParagraphFormat.Bullet.Visible = $false

Write-Host ("Slides: " + $p.Slides.Count)
